$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column style (s="2") from A301 down through the new rows
$ws.Range("A301").Copy($ws.Range("A302:A328"))

$data = @(
    @(302, 44376, 0, 0, 0),
    @(303, 44377, 0, 0, 0),
    @(304, 44378, 1, 1, 48.07692307692308),
    @(305, 44379, 0, 1, 48.07692307692308),
    @(306, 44380, 0, 1, 48.07692307692308),
    @(307, 44381, 0, 1, 48.07692307692308),
    @(308, 44382, 0, 1, 48.07692307692308),
    @(309, 44383, 0, 1, 48.07692307692308),
    @(310, 44384, 0, 1, 48.07692307692308),
    @(311, 44385, 0, 0, 0),
    @(312, 44386, 0, 0, 0),
    @(313, 44387, 0, 0, 0),
    @(314, 44388, 0, 0, 0),
    @(315, 44389, 0, 0, 0),
    @(316, 44390, 0, 0, 0),
    @(317, 44391, 0, 0, 0),
    @(318, 44392, 0, 0, 0),
    @(319, 44393, 0, 0, 0),
    @(320, 44394, 0, 0, 0),
    @(321, 44395, 0, 0, 0),
    @(322, 44396, 0, 0, 0),
    @(323, 44397, 0, 0, 0),
    @(324, 44398, 0, 0, 0),
    @(325, 44399, 0, 0, 0),
    @(326, 44400, 0, 0, 0),
    @(327, 44401, 0, 0, 0),
    @(328, 44402, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}
